# The presentation's design/theme was changed (PowerPoint "Design" gallery):
# the deck's applied colour theme changes from the "Integral" design
# (clrScheme "Red Violet") to the built-in "Office Theme" colour scheme.
#
# PowerPoint stores the document theme in ppt/theme/theme1.xml (referenced by
# the slide master / presentation.xml) and keeps the previously-applied theme
# around as ppt/theme/theme2.xml (referenced only by the notes master) -
# applying a new design only rewrites the 12-slot colour scheme (dk1/lt1/
# dk2/lt2/accent1-6/hlink/folHlink) of the active theme; font scheme and
# format scheme are already shared between the two themes in this deck.
#
# Re-colour the presentation's theme colour scheme (theme1.xml) to the
# standard Office Theme palette via the slide's ThemeColorScheme - this is
# the supported automation surface for recolouring the applied design.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colour scheme (index : role : hex -> decimal BGR-long RGB)
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
